# Generate Report for Handoff
#
# The workbook tracks localization hand-off/hand-back status for two
# source files:
#   FILE_A = 7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.md
#   FILE_B = 15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md
#
# Previously row 2 held FILE_B and row 3 held FILE_A on every sheet (all
# "Handed back: in sync with en-US"). After a new handoff, the two rows
# swap places (row 2 = FILE_A, row 3 = FILE_B) and FILE_B's row reflects
# its fresh "Ready for handoff" status with updated timestamps and a
# staleness warning in the Error Detail column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("A2").Value = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.md"
$ovw.Range("B2").Value = "e2e\7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.md"

$ovw.Range("A3").Value = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md"
$ovw.Range("B3").Value = "e2e\15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md"
$ovw.Range("E3").Value = "Ready for handoff"
$ovw.Range("F3").Value = "Ready for handoff"
$ovw.Range("G3").Value = "2016-08-26 00:47:39"

$ovwLinks = @($ovw.Hyperlinks)
$ovwLinks[0].TextToDisplay = "e2e\7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.md"
$ovwLinks[1].TextToDisplay = "e2e\15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md"

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.md"
$zh.Range("G2").Value = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.bea10080c8689a235a96f1587278e15cad198037.zh-cn.xlf"
$zh.Range("I2").Value = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.md"
$zh.Range("J2").Value = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.bea10080c8689a235a96f1587278e15cad198037.zh-cn.xlf"

$zh.Range("A3").Value = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.420f615339071e6849b4030d2e03f4426b835a8c.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-26 00:47:35"
$zh.Range("I3").Value = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md"
$zh.Range("J3").Value = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.420f615339071e6849b4030d2e03f4426b835a8c.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e82b5fe049d2d0e5badbd7f4eedf563ad7135c25/e2e/15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8b663db51fb028564f4916f94b29823cb07796b/e2e/15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md."

$zhLinks = @($zh.Hyperlinks)
$zhLinks[0].TextToDisplay = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.md"
$zhLinks[1].TextToDisplay = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.md"
$zhLinks[2].TextToDisplay = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md"
$zhLinks[3].TextToDisplay = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md"

$zh.Columns.Item(16).ColumnWidth = 235/6

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.md"
$de.Range("G2").Value = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.bea10080c8689a235a96f1587278e15cad198037.de-de.xlf"
$de.Range("I2").Value = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.md"
$de.Range("J2").Value = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.bea10080c8689a235a96f1587278e15cad198037.de-de.xlf"

$de.Range("A3").Value = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.420f615339071e6849b4030d2e03f4426b835a8c.de-de.xlf"
$de.Range("H3").Value = "2016-08-26 00:47:39"
$de.Range("I3").Value = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md"
$de.Range("J3").Value = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.420f615339071e6849b4030d2e03f4426b835a8c.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e82b5fe049d2d0e5badbd7f4eedf563ad7135c25/e2e/15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8b663db51fb028564f4916f94b29823cb07796b/e2e/15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md."

$deLinks = @($de.Hyperlinks)
$deLinks[0].TextToDisplay = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.md"
$deLinks[1].TextToDisplay = "7c9eefb5-5a68-4fff-b6e2-b1a08a8e7cb7.md"
$deLinks[2].TextToDisplay = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md"
$deLinks[3].TextToDisplay = "15fd3d9f-3bec-4223-9d93-5fd18dccbbf3.md"

$de.Columns.Item(16).ColumnWidth = 235/6
